$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.903.38"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "1.566.32"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.485"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "1.787.57"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.571.05"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "26.896.35"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("E31").Value = "  -3.45%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "1.399.60"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.920"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0164"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.50%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").Value = "1.701.72"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
